# Add a new "Title and Content" slide at the end of the deck
# (becomes slide 12) summarizing the Machine Learning model work.

$p = $ppt.ActivePresentation

# Append a new slide after the last existing slide, using the
# "Title and Content" layout (same layout already used for other
# Title+Content slides in this deck, e.g. slideLayout2.xml).
$lastIndex = $p.Slides.Count
$s = $p.Slides.Add($lastIndex + 1, 2)

# --- Title placeholder -----------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Machine Learning "

# --- Body / content placeholder ---------------------------------------
# Build the paragraphs one at a time via InsertAfter so each stays a
# clean, independent paragraph (keeps the default run-properties intact).
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Data: 291294 observation and 82 features"
$null = $body.InsertAfter("`rTarget: Severity with 4 targets")
$null = $body.InsertAfter("`rPreprocess: handling missing values, remove outliners, feature importance, apply different encoding schemes")
$null = $body.InsertAfter("`rModel: logistic regression, Random forest, and XGB")
$null = $body.InsertAfter("`rResult: Best model is Random forest with accuracy 66% and standard deviation of 0.079")

# Bold the short "label" lead-in of each bullet (matches the source
# deck: bold run followed by a normal run within the same paragraph).
$body.Characters(1, 4).Font.Bold = $true     # "Data"
$body.Characters(42, 8).Font.Bold = $true    # "Target: "
$body.Characters(74, 12).Font.Bold = $true   # "Preprocess: "
$body.Characters(182, 7).Font.Bold = $true   # "Model: "
$body.Characters(233, 6).Font.Bold = $true   # "Result"
